# Scheduled market-price refresh for the Ultima_Profits workbook.
# Re-prices the currentAveragePrice(NQ/HQ) columns (H:N) for the affected
# leve rows on each crafting-class sheet (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR),
# driven by the latest market-board pull. Profit columns (M/N) are derived
# the same way the source sheet already computes them, so they are updated
# alongside the price columns.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 32258378
$ws.Range("I33").Value = 40000370
$ws.Range("J33").Value = 90
$ws.Range("K33").Value = 40000370
$ws.Range("L33").Value = 90
$ws.Range("M33").Value = -40000141
$ws.Range("N33").Value = -548
# Row 53
$ws.Range("H53").Value = 1252.0435
$ws.Range("I53").Value = 2057.818
$ws.Range("J53").Value = 513.4167
$ws.Range("K53").Value = 2057.818
$ws.Range("L53").Value = 513.4167
$ws.Range("M53").Value = -1420.818
$ws.Range("N53").Value = -1787.4167
# Row 129
$ws.Range("H129").Value = 1535.5
$ws.Range("I129").Value = 385
$ws.Range("J129").Value = 1663.3334
$ws.Range("K129").Value = 1155
$ws.Range("L129").Value = 4990.0002
$ws.Range("M129").Value = 3845
$ws.Range("N129").Value = -14990.0002

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 8874.308000000001
$ws.Range("I32").Value = 10247.8125
$ws.Range("K32").Value = 10247.8125
$ws.Range("M32").Value = -9960.8125
# Row 123
$ws.Range("H123").Value = 33427
$ws.Range("J123").Value = 33427
$ws.Range("L123").Value = 33427
$ws.Range("N123").Value = -43227

$ws = $wb.Worksheets.Item("BSM")
# Row 80
$ws.Range("H80").Value = 7428.4644
$ws.Range("I80").Value = 16889.5
$ws.Range("J80").Value = 332.6875
$ws.Range("K80").Value = 16889.5
$ws.Range("L80").Value = 332.6875
$ws.Range("M80").Value = -15891.5
$ws.Range("N80").Value = -2328.6875
# Row 83
$ws.Range("H83").Value = 7428.4644
$ws.Range("I83").Value = 16889.5
$ws.Range("J83").Value = 332.6875
$ws.Range("K83").Value = 84447.5
$ws.Range("L83").Value = 1663.4375
$ws.Range("M83").Value = -79455.5
$ws.Range("N83").Value = -11647.4375
# Row 107
$ws.Range("H107").Value = 4830.143
$ws.Range("I107").Value = 5162.2
$ws.Range("J107").Value = 4000
$ws.Range("K107").Value = 5162.2
$ws.Range("L107").Value = 4000
$ws.Range("M107").Value = -3242.2
$ws.Range("N107").Value = -7840

$ws = $wb.Worksheets.Item("CRP")
# Row 122
$ws.Range("H122").Value = 1407.36
$ws.Range("I122").Value = 1501.6666
$ws.Range("J122").Value = 1164.8572
$ws.Range("K122").Value = 4504.9998
$ws.Range("L122").Value = 3494.5716
$ws.Range("M122").Value = -2054.9998
$ws.Range("N122").Value = -8394.571599999999

$ws = $wb.Worksheets.Item("CUL")
# Row 14
$ws.Range("H14").Value = 923.75
$ws.Range("I14").Value = 923.75
$ws.Range("K14").Value = 2771.25
$ws.Range("M14").Value = -2598.25
# Row 38
$ws.Range("H38").Value = 264.07144
$ws.Range("I38").Value = 416.25
$ws.Range("J38").Value = 203.2
$ws.Range("K38").Value = 1248.75
$ws.Range("L38").Value = 609.5999999999999
$ws.Range("M38").Value = -901.75
$ws.Range("N38").Value = -1303.6
# Row 76
$ws.Range("H76").Value = 3420
$ws.Range("J76").Value = 4125
$ws.Range("L76").Value = 12375
$ws.Range("N76").Value = -13141
# Row 79
$ws.Range("H79").Value = 3420
$ws.Range("J79").Value = 4125
$ws.Range("L79").Value = 12375
$ws.Range("N79").Value = -15027
# Row 99
$ws.Range("H99").Value = 2104.8
$ws.Range("I99").Value = 1762
$ws.Range("J99").Value = 2333.3333
$ws.Range("K99").Value = 5286
$ws.Range("L99").Value = 6999.999899999999
$ws.Range("M99").Value = -3040
$ws.Range("N99").Value = -11491.9999
# Row 104
$ws.Range("H104").Value = 11011.4
$ws.Range("I104").Value = 19999
$ws.Range("J104").Value = 8764.5
$ws.Range("K104").Value = 59997
$ws.Range("L104").Value = 26293.5
$ws.Range("M104").Value = -57376
$ws.Range("N104").Value = -31535.5
# Row 105
$ws.Range("H105").Value = 1785.7142
$ws.Range("I105").Value = 500
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 1500
$ws.Range("L105").Value = 6000
$ws.Range("M105").Value = 1121
$ws.Range("N105").Value = -11242
# Row 107
$ws.Range("H107").Value = 860.43634
$ws.Range("I107").Value = 396.4
$ws.Range("J107").Value = 1417.28
$ws.Range("K107").Value = 1189.2
$ws.Range("L107").Value = 4251.84
$ws.Range("M107").Value = 730.8000000000002
$ws.Range("N107").Value = -8091.84
# Row 108
$ws.Range("H108").Value = 2902.7407
$ws.Range("I108").Value = 998.4
$ws.Range("J108").Value = 4022.9412
$ws.Range("K108").Value = 2995.2
$ws.Range("L108").Value = 12068.8236
$ws.Range("M108").Value = -115.1999999999998
$ws.Range("N108").Value = -17828.8236
# Row 111
$ws.Range("H111").Value = 0
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("L111").ClearContents()
$ws.Range("M111").ClearContents()
$ws.Range("N111").Value = 0
# Row 112
$ws.Range("H112").Value = 7134.5713
$ws.Range("J112").Value = 6648
$ws.Range("L112").Value = 19944
$ws.Range("N112").Value = -22160
# Row 113
$ws.Range("H113").Value = 1235.7778
$ws.Range("J113").Value = 1687
$ws.Range("L113").Value = 5061
$ws.Range("N113").Value = -9401
# Row 114
$ws.Range("H114").Value = 1283.4286
$ws.Range("I114").Value = 314
$ws.Range("J114").Value = 3028.4
$ws.Range("K114").Value = 942
$ws.Range("L114").Value = 9085.200000000001
$ws.Range("M114").Value = 2312
$ws.Range("N114").Value = -15593.2
# Row 116
$ws.Range("H116").Value = 2124.8667
$ws.Range("J116").Value = 3685.1428
$ws.Range("L116").Value = 11055.4284
$ws.Range("N116").Value = -17939.4284
# Row 117
$ws.Range("H117").Value = 1137.5294
$ws.Range("I117").Value = 442.33334
$ws.Range("J117").Value = 1286.5
$ws.Range("K117").Value = 1327.00002
$ws.Range("L117").Value = 3859.5
$ws.Range("M117").Value = 2114.99998
$ws.Range("N117").Value = -10743.5
# Row 118
$ws.Range("H118").Value = 1882.7273
$ws.Range("I118").Value = 1525
$ws.Range("J118").Value = 1918.5
$ws.Range("K118").Value = 4575
$ws.Range("L118").Value = 5755.5
$ws.Range("M118").Value = -3332
$ws.Range("N118").Value = -8241.5
# Row 119
$ws.Range("H119").Value = 15021.75
$ws.Range("I119").Value = 10043.5
$ws.Range("K119").Value = 30130.5
$ws.Range("M119").Value = -25292.5
# Row 120
$ws.Range("H120").Value = 12500
$ws.Range("I120").Value = 5000
$ws.Range("J120").Value = 20000
$ws.Range("K120").Value = 15000
$ws.Range("L120").Value = 60000
$ws.Range("M120").Value = -10162
$ws.Range("N120").Value = -69676
# Row 122
$ws.Range("H122").Value = 879.26086
$ws.Range("I122").Value = 1261.6364
$ws.Range("J122").Value = 528.75
$ws.Range("K122").Value = 11354.7276
$ws.Range("L122").Value = 4758.75
$ws.Range("M122").Value = -8904.7276
$ws.Range("N122").Value = -9658.75

$ws = $wb.Worksheets.Item("GSM")
# Row 48
$ws.Range("H48").Value = 89338
$ws.Range("J48").Value = 19007
$ws.Range("L48").Value = 19007
$ws.Range("N48").Value = -19977
# Row 123
$ws.Range("H123").Value = 23275.715
$ws.Range("J123").Value = 23275.715
$ws.Range("L123").Value = 23275.715
$ws.Range("N123").Value = -28175.715
# Row 126
$ws.Range("H126").Value = 5646
$ws.Range("I126").Value = 4000
$ws.Range("J126").Value = 5783.1665
$ws.Range("K126").Value = 12000
$ws.Range("L126").Value = 17349.4995
$ws.Range("M126").Value = -9530
$ws.Range("N126").Value = -22289.4995

$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 7047674.5
$ws.Range("I132").Value = 3249.5
$ws.Range("J132").Value = 26327154
$ws.Range("K132").Value = 9748.5
$ws.Range("L132").Value = 78981462
$ws.Range("M132").Value = -7218.5
$ws.Range("N132").Value = -78986522
# Row 136
$ws.Range("H136").Value = 32613804
$ws.Range("I136").Value = 46876844
$ws.Range("J136").Value = 12567.857
$ws.Range("K136").Value = 140630532
$ws.Range("L136").Value = 37703.571
$ws.Range("M136").Value = -140627982
$ws.Range("N136").Value = -42803.571

$ws = $wb.Worksheets.Item("WVR")
# Row 28
$ws.Range("H28").Value = 3128141.5
$ws.Range("J28").Value = 3128141.5
$ws.Range("L28").Value = 3128141.5
$ws.Range("N28").Value = -3128837.5
# Row 74
$ws.Range("H74").Value = 9666.833000000001
$ws.Range("J74").Value = 9666.833000000001
$ws.Range("L74").Value = 9666.833000000001
$ws.Range("N74").Value = -11538.833
# Row 77
$ws.Range("H77").Value = 9666.833000000001
$ws.Range("J77").Value = 9666.833000000001
$ws.Range("L77").Value = 29000.499
$ws.Range("N77").Value = -38360.499
# Row 123
$ws.Range("H123").Value = 42772.145
$ws.Range("J123").Value = 42772.145
$ws.Range("L123").Value = 42772.145
$ws.Range("N123").Value = -52572.145
# Row 126
$ws.Range("H126").Value = 3084.9333
$ws.Range("I126").Value = 1932.2354
$ws.Range("J126").Value = 4592.3076
$ws.Range("K126").Value = 5796.706200000001
$ws.Range("L126").Value = 13776.9228
$ws.Range("M126").Value = -3326.706200000001
$ws.Range("N126").Value = -18716.9228
